$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# ------------------------------------------------------------------
# 1) Two fixtures played on the same day (match ids 6627737 and
#    6627736, both "Bulgaria First League" games on 2023-06-06) have
#    their full row content (every column except the running "id" in
#    column A) swapped between the two rows.
# ------------------------------------------------------------------
$rowA = -1
$rowB = -1
for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    if ($bVal -eq 6627737) { $rowA = $r }
    if ($bVal -eq 6627736) { $rowB = $r }
    if ($rowA -gt 0 -and $rowB -gt 0) { break }
}

if ($rowA -gt 0 -and $rowB -gt 0) {
    $rangeA = $ws.Range("B" + $rowA + ":AC" + $rowA)
    $rangeB = $ws.Range("B" + $rowB + ":AC" + $rowB)
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

# ------------------------------------------------------------------
# 2) Remove the fixture row whose match id (column B) is 6978388
#    (FC Hebar Pazardzhik vs Etar 1924 Veliko Tarnovo on
#    2024-04-06). All following rows shift up one position.
# ------------------------------------------------------------------
$targetRow = -1
for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    if ($bVal -eq 6978388) {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    # Decrement the running id (column A) for every row below the one
    # that is about to be removed, so ids stay contiguous after the
    # row shift caused by the delete.
    for ($r = $targetRow + 1; $r -le $lastRow; $r++) {
        $idCell = $ws.Cells.Item($r, 1)
        $idVal = $idCell.Value()
        if ($idVal -ne $null) {
            $idCell.Value = $idVal - 1
        }
    }

    $ws.Rows.Item($targetRow).Delete()
}
